$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve the original (empty) F1 cell; the COM round-trip otherwise
# resolves the empty shared-string reference to the first shared string.
$ws.Cells.Item(1, 6).Value = ""

# Add "NA" values under the duplicate_image_filename column (column E) for rows 2 through 21
for ($row = 2; $row -le 21; $row++) {
    $ws.Cells.Item($row, 5).Value = "NA"
}
